$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the F:V data between rows 21 and 22 (columns A-E unchanged) ---
# The same fixture data for 06-08/08/2023 had its two matches (Randers FC vs
# Nordsjaelland, and Vejle vs Midtjylland) recorded under the wrong rows; this
# swaps the match details (F:V) back between row 21 and row 22.
$cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")
foreach ($col in $cols) {
    $v21 = $ws.Range($col + "21").Value2
    $v22 = $ws.Range($col + "22").Value2
    $ws.Range($col + "21").Value2 = $v22
    $ws.Range($col + "22").Value2 = $v21
}

# --- Append 5 new match rows (99-103), copying the style from row 98 first ---
$ws.Range("A98:V98").Copy()
$ws.Range("A99:V103").PasteSpecial(-4122)

# Row 99 (Indice=98)
$ws.Range("A99").Value2 = 98
$ws.Range("B99").Value2 = 'denmark'
$ws.Range("C99").Value2 = 'superliga'
$ws.Range("D99").Value2 = '2023-2024'
$ws.Range("E99").Value2 = 45263.58333333334
$ws.Range("F99").Value2 = 'Lyngby'
$ws.Range("G99").Value2 = 2
$ws.Range("H99").Value2 = 'Silkeborg'
$ws.Range("I99").Value2 = 0
$ws.Range("J99").Value2 = 3.52
$ws.Range("K99").Value2 = '27/11/2023 19:12'
$ws.Range("L99").Value2 = 3.15
$ws.Range("M99").Value2 = '03/12/2023 13:55'
$ws.Range("N99").Value2 = 3.67
$ws.Range("O99").Value2 = '27/11/2023 19:12'
$ws.Range("P99").Value2 = 3.54
$ws.Range("Q99").Value2 = '03/12/2023 13:58'
$ws.Range("R99").Value2 = 1.98
$ws.Range("S99").Value2 = '27/11/2023 19:12'
$ws.Range("T99").Value2 = 2.31
$ws.Range("U99").Value2 = '03/12/2023 13:55'
$ws.Range("V99").Value2 = 'https://www.betexplorer.com/football/denmark/superliga/lyngby-silkeborg/6DXBMDG8/'

# Row 100 (Indice=99)
$ws.Range("A100").Value2 = 99
$ws.Range("B100").Value2 = 'denmark'
$ws.Range("C100").Value2 = 'superliga'
$ws.Range("D100").Value2 = '2023-2024'
$ws.Range("E100").Value2 = 45263.58333333334
$ws.Range("F100").Value2 = 'Odense'
$ws.Range("G100").Value2 = 1
$ws.Range("H100").Value2 = 'Nordsjaelland'
$ws.Range("I100").Value2 = 1
$ws.Range("J100").Value2 = 4.04
$ws.Range("K100").Value2 = '27/11/2023 09:49'
$ws.Range("L100").Value2 = 4.09
$ws.Range("M100").Value2 = '03/12/2023 13:56'
$ws.Range("N100").Value2 = 3.84
$ws.Range("O100").Value2 = '27/11/2023 09:49'
$ws.Range("P100").Value2 = 3.65
$ws.Range("Q100").Value2 = '03/12/2023 13:56'
$ws.Range("R100").Value2 = 1.79
$ws.Range("S100").Value2 = '27/11/2023 09:49'
$ws.Range("T100").Value2 = 1.93
$ws.Range("U100").Value2 = '03/12/2023 13:56'
$ws.Range("V100").Value2 = 'https://www.betexplorer.com/football/denmark/superliga/odense-nordsjaelland/6m2GJB0R/'

# Row 101 (Indice=100)
$ws.Range("A101").Value2 = 100
$ws.Range("B101").Value2 = 'denmark'
$ws.Range("C101").Value2 = 'superliga'
$ws.Range("D101").Value2 = '2023-2024'
$ws.Range("E101").Value2 = 45263.66666666666
$ws.Range("F101").Value2 = 'Brondby'
$ws.Range("G101").Value2 = 4
$ws.Range("H101").Value2 = 'Hvidovre IF'
$ws.Range("I101").Value2 = 0
$ws.Range("J101").Value2 = 1.26
$ws.Range("K101").Value2 = '26/11/2023 18:13'
$ws.Range("L101").Value2 = 1.2
$ws.Range("M101").Value2 = '03/12/2023 15:57'
$ws.Range("N101").Value2 = 6.1
$ws.Range("O101").Value2 = '26/11/2023 18:13'
$ws.Range("P101").Value2 = 7.28
$ws.Range("Q101").Value2 = '03/12/2023 15:57'
$ws.Range("R101").Value2 = 10.74
$ws.Range("S101").Value2 = '26/11/2023 18:13'
$ws.Range("T101").Value2 = 13.91
$ws.Range("U101").Value2 = '03/12/2023 15:57'
$ws.Range("V101").Value2 = 'https://www.betexplorer.com/football/denmark/superliga/brondby-hvidovre-if/lx3RE1c2/'

# Row 102 (Indice=101)
$ws.Range("A102").Value2 = 101
$ws.Range("B102").Value2 = 'denmark'
$ws.Range("C102").Value2 = 'superliga'
$ws.Range("D102").Value2 = '2023-2024'
$ws.Range("E102").Value2 = 45263.75
$ws.Range("F102").Value2 = 'FC Copenhagen'
$ws.Range("G102").Value2 = 1
$ws.Range("H102").Value2 = 'Aarhus'
$ws.Range("I102").Value2 = 2
$ws.Range("J102").Value2 = 1.71
$ws.Range("K102").Value2 = '26/11/2023 14:13'
$ws.Range("L102").Value2 = 1.67
$ws.Range("M102").Value2 = '03/12/2023 17:56'
$ws.Range("N102").Value2 = 3.83
$ws.Range("O102").Value2 = '26/11/2023 14:13'
$ws.Range("P102").Value2 = 3.69
$ws.Range("Q102").Value2 = '03/12/2023 17:57'
$ws.Range("R102").Value2 = 4.53
$ws.Range("S102").Value2 = '26/11/2023 14:13'
$ws.Range("T102").Value2 = 5.97
$ws.Range("U102").Value2 = '03/12/2023 17:57'
$ws.Range("V102").Value2 = 'https://www.betexplorer.com/football/denmark/superliga/fc-copenhagen-aarhus/vRVJKioL/'

# Row 103 (Indice=102)
$ws.Range("A103").Value2 = 102
$ws.Range("B103").Value2 = 'denmark'
$ws.Range("C103").Value2 = 'superliga'
$ws.Range("D103").Value2 = '2023-2024'
$ws.Range("E103").Value2 = 45264.79166666666
$ws.Range("F103").Value2 = 'Midtjylland'
$ws.Range("G103").Value2 = 5
$ws.Range("H103").Value2 = 'Viborg'
$ws.Range("I103").Value2 = 1
$ws.Range("J103").Value2 = 1.53
$ws.Range("K103").Value2 = '27/11/2023 19:12'
$ws.Range("L103").Value2 = 1.54
$ws.Range("M103").Value2 = '04/12/2023 18:58'
$ws.Range("N103").Value2 = 4.26
$ws.Range("O103").Value2 = '27/11/2023 19:12'
$ws.Range("P103").Value2 = 4.35
$ws.Range("Q103").Value2 = '04/12/2023 18:56'
$ws.Range("R103").Value2 = 5.57
$ws.Range("S103").Value2 = '27/11/2023 19:12'
$ws.Range("T103").Value2 = 6.36
$ws.Range("U103").Value2 = '04/12/2023 18:58'
$ws.Range("V103").Value2 = 'https://www.betexplorer.com/football/denmark/superliga/midtjylland-viborg/hAT7Ng12/'
